# Logic change for Logged in User
# Insert a new "CLICK LoginURL" step before the Uname entry step, and
# rename the Uname1/Password1/LoginButton1 objects to Uname/Password/LoginButton.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC31_Verify_login")

# Insert a new row above the current row 4 (ENTERTEXT / Uname1 ...),
# shifting the remaining steps down by one.
$ws.Rows.Item(4).Insert()

# New step: CLICK the LoginURL object (CSS) before entering credentials.
$ws.Range("B4").Value = "CLICK"
$ws.Range("C4").Value = "LoginURL"
$ws.Range("D4").Value = "CSS"

# Give the newly inserted row the same look (thin border all round) as
# the rest of the data rows, by copying the border from the row above.
$ws.Range("A3:E3").Copy()
$ws.Range("A4:E4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The remaining steps shifted down to rows 5-12; rename the Uname1 /
# Password1 / LoginButton1 objects to Uname / Password / LoginButton.
$ws.Range("C5").Value = "Uname"
$ws.Range("E5").Value = "Uname"

$ws.Range("C6").Value = "Password"
$ws.Range("E6").Value = "Password"

$ws.Range("C7").Value = "LoginButton"

# Match the saved selection left behind by the authoring session.
$ws.Activate()
$ws.Range("C5:C7").Select()
